$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Alternativo 2" step text: passo 13 -> passo 11
$ws.Range("B20").Value = "Alternativo 2 [Não confirma compra] (passo 11)"

# Update the sub-step text: 13.1 -> 11.1
$ws.Range("C20").Value = "11.1 Não confirma compra"

# Update the "Regressa a" text: Regressa a 9 -> Regressa a 1
$ws.Range("D21").Value = "Regressa a 1"

# Update view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("D21").Select()
